$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Update title (B1, shared string reused)
$ws.Range("B1").Value = "BAJA SAE 2017-2018 Electronics Bill of Materials"

# 2. Move the old "Total" label (row 15, D column only) down to D17, keep its
#    format. Avoid touching E15/E17 with a copied formula - re-creating a
#    formula value in a cell that previously held a different cached formula
#    can leave the recalculated value stale, so E17 is built fresh later and
#    only has a style pasted onto it as the very last step.
$ws.Range("D15").Copy()
$ws.Range("D17").PasteSpecial(-4122)
$ws.Range("D15:E15").ClearContents()

# 3. Prime rows 15 and 16 with the same per-column formatting as the other
#    data rows (columns B-H), by copying the format of row 3.
$ws.Range("B3:H3").Copy($ws.Range("B15:H15"))
$ws.Range("B3:H3").Copy($ws.Range("B16:H16"))

# 4. New column I header + data formatting: copy format from the currency
#    header style (C2) for I2, and from a plain data cell (D3) for I3:I16.
$ws.Range("C2").Copy($ws.Range("I2"))
$ws.Range("D3").Copy($ws.Range("I3"))
$ws.Range("D3").Copy($ws.Range("I4"))
$ws.Range("D3").Copy($ws.Range("I5"))
$ws.Range("D3").Copy($ws.Range("I6"))
$ws.Range("D3").Copy($ws.Range("I7"))
$ws.Range("D3").Copy($ws.Range("I8"))
$ws.Range("D3").Copy($ws.Range("I9"))
$ws.Range("D3").Copy($ws.Range("I10"))
$ws.Range("D3").Copy($ws.Range("I11"))
$ws.Range("D3").Copy($ws.Range("I12"))
$ws.Range("D3").Copy($ws.Range("I13"))
$ws.Range("D3").Copy($ws.Range("I14"))
$ws.Range("D3").Copy($ws.Range("I15"))
$ws.Range("D3").Copy($ws.Range("I16"))

# 5. Column I values
$ws.Range("I2").Value = "Power Source"
$ws.Range("I3").Value = "Arduino 5V"
$ws.Range("I4").Value = "Arduino 5V"
$ws.Range("I5").Value = "External 5V"
$ws.Range("I6").Value = "External 5V"
$ws.Range("I7").Value = "External 5V"
$ws.Range("I8").Value = "Arduino 5V"
$ws.Range("I9").Value = "External 5V"
$ws.Range("I10").Value = "None"
$ws.Range("I11").Value = "None"
$ws.Range("I12").Value = "Arduino 5V"
$ws.Range("I13").Value = "External 5V"
$ws.Range("I14").Value = "None"
$ws.Range("I15").Value = "External 9V"
$ws.Range("I16").Value = "External 12V"

# 6. Update quantities for rows 10 and 11 (On/Off Switch, Arcade Button): 1 -> 2
$ws.Range("D10").Value = 2
$ws.Range("D11").Value = 2

# 7. Fill in new row 15 - Arduino Mega
$ws.Range("B15").Value = "Arduino Mega"
$ws.Range("C15").Value = 38.5
$ws.Range("D15").Value = 1
$ws.Range("F15").Value = "arduino"
$ws.Hyperlinks.Add($ws.Range("G15"), "https://store.arduino.cc/usa/arduino-mega-2560-rev3", "", "", "https://store.arduino.cc/usa/arduino-mega-2560-rev3")
$ws.Range("G15").Style = "Hyperlink"
$ws.Range("H15").Value = "Control all input and output for components"

# 8. Fill in new row 16 - Computer Fan
$ws.Range("B16").Value = "Computer Fan"
$ws.Range("C16").Value = 23.39
$ws.Range("D16").Value = 1
$ws.Range("F16").Value = "PC Hub"
$ws.Hyperlinks.Add($ws.Range("G16"), "http://www.pchub.com/uph/laptop/656-79749-22693/Cooler-Master-MGT8012ZR-W25-Server-Square-Fan.html", "", "", "http://www.pchub.com/uph/laptop/656-79749-22693/Cooler-Master-MGT8012ZR-W25-Server-Square-Fan.html")
$ws.Range("G16").Style = "Hyperlink"
$ws.Range("H16").Value = "Cool components"

# 9. Formulas: extend the shared E7:E14 pattern down through E16.
$ws.Range("E7:E16").FormulaR1C1 = "=RC[-2]*RC[-1]"

# 10. New Total row 17 - value/formula first, format pasted on last so the
#     recalculated SUM isn't left stale.
$ws.Range("D17").Value = "Total"
$ws.Range("E17").Formula = "=SUM(E3:E16)"
$ws.Range("C2").Copy()
$ws.Range("E17").PasteSpecial(-4122)

# 11. Column widths
$ws.Columns.Item(8).ColumnWidth = 24
$ws.Columns.Item(9).ColumnWidth = 14.5703125

# 12. Update selection
$ws.Range("J2").Select()
